$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all updated cells keep their original text formatting (no numeric auto-conversion)

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.942.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.21%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.743.61"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.37%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.50%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.15"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.24%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.109"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.98%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.24%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.67"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.90%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.382"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.13%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.227.71"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.40%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.49"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.05%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.606.28"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.58%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000150"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.12%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.745.07"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.33%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.07"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.49%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.80"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.84%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.66"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.88%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.57"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.81%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.15%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.521"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -6.66%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.45"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.22%  "

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.81%  "

# Row 26
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.06%  "

# Row 27
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.42"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.76%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0910"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.03%  "

# Row 29
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.97"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.32%  "

# Row 30
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.20"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.34%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.35"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +8.53%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.61"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.72%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.90"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.49%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.04"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.68%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.48"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.63%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.80"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.29%  "

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "352.34"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.59%  "

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "SuiNetwork"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.985"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.58%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.24"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.00%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.13"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.06%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.55"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.41%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.10"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.28%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.23"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.46%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0583"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.63%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.90"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.68%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.626"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.74%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.100"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.32%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0248"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.65%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.00%  "

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.05"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.06%  "
